# Applies the value updates to the Leve profit-tracking tables (columns H:N)
# across all eight crafting-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Source data cells (A:G) are untouched; only cached currentAveragePrice /
# LevePrice* / LeveProfit* figures are refreshed to the latest scrape.

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 250.5
$ws.Range("I12").Value = 250.5
$ws.Range("K12").Value = 250.5
$ws.Range("M12").Value = -80.5
# Row 17
$ws.Range("H17").Value = 1812.5
$ws.Range("J17").Value = 1812.5
$ws.Range("L17").Value = 5437.5
$ws.Range("N17").Value = -5773.5
# Row 40
$ws.Range("H40").Value = 3572.3333
$ws.Range("I40").Value = 3089.5
$ws.Range("J40").Value = 4538
$ws.Range("K40").Value = 3089.5
$ws.Range("L40").Value = 4538
$ws.Range("M40").Value = -2914.5
$ws.Range("N40").Value = -4888
# Row 74
$ws.Range("H74").Value = 2367.6667
$ws.Range("I74").Value = 2367.6667
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2367.6667
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1431.6667
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 2367.6667
$ws.Range("I77").Value = 2367.6667
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 11838.3335
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -7158.333500000001
$ws.Range("N77").ClearContents()
# Row 100
$ws.Range("H100").Value = 1365.8334
$ws.Range("I100").Value = 1365.8334
$ws.Range("K100").Value = 1365.8334
$ws.Range("M100").Value = -824.8334
# Row 113
$ws.Range("H113").Value = 166667660
$ws.Range("J113").Value = 200001500
$ws.Range("L113").Value = 200001500
$ws.Range("N113").Value = -200008008
# Row 124
$ws.Range("H124").Value = 90390
$ws.Range("J124").Value = 90390
$ws.Range("L124").Value = 90390
$ws.Range("N124").Value = -100210
# Row 127
$ws.Range("H127").Value = 15000
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
# Row 129
$ws.Range("H129").Value = 1457.3
$ws.Range("I129").Value = 946.875
$ws.Range("J129").Value = 3499
$ws.Range("K129").Value = 2840.625
$ws.Range("L129").Value = 10497
$ws.Range("M129").Value = 2159.375
$ws.Range("N129").Value = -20497
# Row 132
$ws.Range("H132").Value = 1057.0714
$ws.Range("I132").Value = 929.55554
$ws.Range("K132").Value = 2788.66662
$ws.Range("M132").Value = -258.66662
# Row 137
$ws.Range("H137").Value = 3673.077
$ws.Range("I137").Value = 2506.5898
$ws.Range("K137").Value = 7519.769400000001
$ws.Range("M137").Value = -4969.769400000001

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 737.6
$ws.Range("I2").Value = 615.8125
$ws.Range("J2").Value = 1224.75
$ws.Range("K2").Value = 615.8125
$ws.Range("L2").Value = 1224.75
$ws.Range("M2").Value = -502.8125
$ws.Range("N2").Value = -1450.75
# Row 32
$ws.Range("H32").Value = 13520231
$ws.Range("I32").Value = 19233798
$ws.Range("K32").Value = 19233798
$ws.Range("M32").Value = -19233511
# Row 37
$ws.Range("H37").Value = 49285.715
$ws.Range("I37").Value = 48000
$ws.Range("K37").Value = 48000
$ws.Range("M37").Value = -47727
# Row 46
$ws.Range("H46").Value = 15500
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 15500
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 15500
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -16138
# Row 63
$ws.Range("H63").Value = 8101.222
$ws.Range("I63").Value = 2399
$ws.Range("K63").Value = 2399
$ws.Range("M63").Value = -1713
# Row 66
$ws.Range("H66").Value = 8101.222
$ws.Range("I66").Value = 2399
$ws.Range("K66").Value = 11995
$ws.Range("M66").Value = -8563
# Row 74
$ws.Range("H74").Value = 15300161
$ws.Range("I74").Value = 25003008
$ws.Range("K74").Value = 25003008
$ws.Range("M74").Value = -25002134
# Row 77
$ws.Range("H77").Value = 15300161
$ws.Range("I77").Value = 25003008
$ws.Range("K77").Value = 125015040
$ws.Range("M77").Value = -125010672
# Row 116
$ws.Range("H116").Value = 737.6
$ws.Range("I116").Value = 615.8125
$ws.Range("J116").Value = 1224.75
$ws.Range("K116").Value = 615.8125
$ws.Range("L116").Value = 1224.75
$ws.Range("M116").Value = 1678.1875
$ws.Range("N116").Value = -5812.75
# Row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 737.6
$ws.Range("I3").Value = 615.8125
$ws.Range("J3").Value = 1224.75
$ws.Range("K3").Value = 615.8125
$ws.Range("L3").Value = 1224.75
$ws.Range("M3").Value = -501.8125
$ws.Range("N3").Value = -1452.75

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 5008
$ws.Range("I17").Value = 5008
$ws.Range("K17").Value = 5008
$ws.Range("M17").Value = -4834
# Row 31
$ws.Range("H31").Value = 757651.2
$ws.Range("I31").Value = 3162.3572
$ws.Range("K31").Value = 3162.3572
$ws.Range("M31").Value = -2867.3572
# Row 33
$ws.Range("H33").Value = 4018.6
# Row 34
$ws.Range("H34").Value = 757651.2
$ws.Range("I34").Value = 3162.3572
$ws.Range("K34").Value = 3162.3572
$ws.Range("M34").Value = -2960.3572
# Row 132
$ws.Range("H132").Value = 2928.8823
$ws.Range("I132").Value = 2878
$ws.Range("K132").Value = 8634
$ws.Range("M132").Value = -6104
# Row 134
$ws.Range("H134").Value = 2640.5
$ws.Range("I134").Value = 2092.037
$ws.Range("K134").Value = 6276.110999999999
$ws.Range("M134").Value = -3741.110999999999
# Row 140
$ws.Range("H140").Value = 72050
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 42
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
# Row 105
$ws.Range("H105").Value = 9400
$ws.Range("J105").Value = 9400
$ws.Range("L105").Value = 28200
$ws.Range("N105").Value = -33442
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
# Row 122
$ws.Range("H122").Value = 2230.4517
$ws.Range("I122").Value = 491.5
$ws.Range("J122").Value = 3328.7368
$ws.Range("K122").Value = 4423.5
$ws.Range("L122").Value = 29958.6312
$ws.Range("M122").Value = -1973.5
$ws.Range("N122").Value = -34858.6312
# Row 131
$ws.Range("H131").Value = 9785.388999999999
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 303.33334
$ws.Range("I2").Value = 81
$ws.Range("J2").Value = 366.85715
$ws.Range("K2").Value = 81
$ws.Range("L2").Value = 366.85715
$ws.Range("M2").Value = 32
$ws.Range("N2").Value = -592.85715
# Row 41
$ws.Range("H41").Value = 3500
$ws.Range("I41").Value = 3500
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 3500
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -3145
$ws.Range("N41").ClearContents()
# Row 93
$ws.Range("H93").Value = 64437.75
$ws.Range("J93").Value = 64437.75
$ws.Range("L93").Value = 64437.75
$ws.Range("N93").Value = -68181.75
# Row 102
$ws.Range("H102").Value = 5660.3335
$ws.Range("I102").Value = 6737.5
$ws.Range("K102").Value = 6737.5
$ws.Range("M102").Value = -5115.5
# Row 132
$ws.Range("H132").Value = 19236382
$ws.Range("I132").Value = 28574450
$ws.Range("K132").Value = 85723350
$ws.Range("M132").Value = -85720820

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 36
$ws.Range("H36").Value = 79000
$ws.Range("J36").Value = 79000
$ws.Range("L36").Value = 79000
$ws.Range("N36").Value = -80124
# Row 40
$ws.Range("H40").Value = 5425.3335
$ws.Range("I40").Value = 5034.6665
$ws.Range("K40").Value = 5034.6665
$ws.Range("M40").Value = -4898.6665
# Row 82
$ws.Range("H82").Value = 1100
$ws.Range("I82").Value = 990
$ws.Range("K82").Value = 990
$ws.Range("M82").Value = -629
# Row 85
$ws.Range("H85").Value = 1100
$ws.Range("I85").Value = 990
$ws.Range("K85").Value = 990
$ws.Range("M85").Value = 258
# Row 93
$ws.Range("H93").Value = 142858600
$ws.Range("I93").Value = 250001330
$ws.Range("J93").Value = 1590.6666
$ws.Range("K93").Value = 250001330
$ws.Range("L93").Value = 1590.6666
$ws.Range("M93").Value = -250000082
$ws.Range("N93").Value = -4086.6666
# Row 100
$ws.Range("H100").Value = 4122.7
$ws.Range("I100").Value = 3570.75
$ws.Range("K100").Value = 3570.75
$ws.Range("M100").Value = -3029.75
# Row 122
$ws.Range("H122").Value = 5845.885
$ws.Range("I122").Value = 5388.4116
$ws.Range("K122").Value = 16165.2348
$ws.Range("M122").Value = -13715.2348
# Row 136
$ws.Range("H136").Value = 120425.5
$ws.Range("J136").Value = 221835
$ws.Range("L136").Value = 665505
$ws.Range("N136").Value = -670605

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 33353352
$ws.Range("I6").Value = 30027.5
$ws.Range("J6").Value = 100000000
$ws.Range("K6").Value = 30027.5
$ws.Range("L6").Value = 100000000
$ws.Range("M6").Value = -29912.5
$ws.Range("N6").Value = -100000230
# Row 122
$ws.Range("H122").Value = 2442.3914
$ws.Range("I122").Value = 2359.4443
$ws.Range("J122").Value = 2741
$ws.Range("K122").Value = 7078.3329
$ws.Range("L122").Value = 8223
$ws.Range("M122").Value = -4628.3329
$ws.Range("N122").Value = -13123
# Row 123
$ws.Range("H123").Value = 75050
$ws.Range("J123").Value = 75050
$ws.Range("L123").Value = 75050
$ws.Range("N123").Value = -84850
# Row 128
$ws.Range("H128").Value = 156000
$ws.Range("J128").Value = 156000
$ws.Range("L128").Value = 156000
$ws.Range("N128").Value = -165960
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
